# Generate Report for Handoff
# This script applies the "handback" localization status update:
#   - a new file (8e483ecc-d59f-46e1-a834-0050fbf7442f.png) replaces
#     the old source entry at row 2,
#   - a second new file (e733822b-896f-4ded-9014-14571e198879.png) is
#     inserted as a new row,
#   - the original markdown file (now renamed/rehashed to
#     ed1e4cf4-08eb-484d-82d0-7b25239db4f7.md) is re-added as its own row,
#   - the .localization-config row is pushed down to stay last.
# The same structural change is applied on the Overview sheet and on both
# locale detail sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$commitSha = "7b0a4d5f6171c6fd2fe40a1775311e128edfa23f"
$zhHandoffSha = "4d6819da8ce4f34916d08b5e5e91feb6700d7633"
$deHandoffSha = "c5906712996225ffacb82dc297fe259496da2ada"

$srcBase = "https://github.com/OpenLocalizationTest/oltest/blob/$commitSha"
$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhHandoffSha/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deHandoffSha/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$pngFile1 = "8e483ecc-d59f-46e1-a834-0050fbf7442f.png"
$pngFile2 = "e733822b-896f-4ded-9014-14571e198879.png"
$mdFile = "ed1e4cf4-08eb-484d-82d0-7b25239db4f7.md"
$cfgFile = ".localization-config"

$zhXlf1 = "4205b8f53be0b9d2bc57f21de5d4eb356190701c.png"
$zhXlf2 = "bdad5253a96b46ef5f9b98fe44eab92921264d4b.png"
$zhXlf3 = "ed1e4cf4-08eb-484d-82d0-7b25239db4f7.f465581ac5efeb4838bbc130851f1576badaa4f1.zh-cn.xlf"

$deXlf1 = "4205b8f53be0b9d2bc57f21de5d4eb356190701c.png"
$deXlf2 = "bdad5253a96b46ef5f9b98fe44eab92921264d4b.png"
$deXlf3 = "ed1e4cf4-08eb-484d-82d0-7b25239db4f7.f465581ac5efeb4838bbc130851f1576badaa4f1.de-de.xlf"

$zhDatetime = "2016-03-09 05:19:53"
$deDatetime = "2016-03-09 05:19:56"
$epochDatetime = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Push the ".localization-config" row (currently row 3) down to row 5 by
# inserting two blank rows above it; this copies the existing row
# formatting (hyperlink-style column A) down onto the new rows.
$ws1.Rows("3:4").Insert()

$ws1.Range("A2").Value = $pngFile1
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("A3").Value = $pngFile2
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("A4").Value = $mdFile
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

# Row 5 (old row 3) already carries ".localization-config" / "Not to be
# localized" values; nothing else needs to change there.

# Rebuild every hyperlink on the sheet so each "display" ref matches the
# cell's new position (the insert above does not retarget hyperlink refs).
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "$srcBase/e2e/$pngFile1", "", "", $pngFile1)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$srcBase/e2e/$pngFile2", "", "", $pngFile2)
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$srcBase/e2e/$mdFile", "", "", $mdFile)
$ws1.Hyperlinks.Add($ws1.Range("A5"), "$srcBase/$cfgFile", "", "", $cfgFile)

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows("3:4").Insert()

$ws2.Range("A2").Value = $pngFile1
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = $zhXlf1
$ws2.Range("D2").Value = $zhDatetime
$ws2.Range("G2").Value = $epochDatetime
$ws2.Range("H2").Value = "IsDependency"
$ws2.Range("I2").Value = "e2e\$mdFile"

$ws2.Range("A3").Value = $pngFile2
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = $zhXlf2
$ws2.Range("D3").Value = $zhDatetime
$ws2.Range("G3").Value = $epochDatetime
$ws2.Range("H3").Value = "IsDependency"
$ws2.Range("I3").Value = "e2e\$mdFile"

$ws2.Range("A4").Value = $mdFile
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = $zhXlf3
$ws2.Range("D4").Value = $zhDatetime
$ws2.Range("G4").Value = $epochDatetime
$ws2.Range("H4").Value = "Include"

# Row 5 (old row 3) keeps its ".localization-config" / "Not to be
# localized" / "Ignored" values.

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "$srcBase/e2e/$pngFile1", "", "", $pngFile1)
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$zhHandoffBase/$zhXlf1", "", "", $zhXlf1)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$srcBase/e2e/$pngFile2", "", "", $pngFile2)
$ws2.Hyperlinks.Add($ws2.Range("C3"), "$zhHandoffBase/$zhXlf2", "", "", $zhXlf2)
$ws2.Hyperlinks.Add($ws2.Range("A4"), "$srcBase/e2e/$mdFile", "", "", $mdFile)
$ws2.Hyperlinks.Add($ws2.Range("C4"), "$zhHandoffBase/$zhXlf3", "", "", $zhXlf3)
$ws2.Hyperlinks.Add($ws2.Range("A5"), "$srcBase/$cfgFile", "", "", $cfgFile)

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Rows("3:4").Insert()

$ws3.Range("A2").Value = $pngFile1
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = $deXlf1
$ws3.Range("D2").Value = $deDatetime
$ws3.Range("G2").Value = $epochDatetime
$ws3.Range("H2").Value = "IsDependency"
$ws3.Range("I2").Value = "e2e\$mdFile"

$ws3.Range("A3").Value = $pngFile2
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = $deXlf2
$ws3.Range("D3").Value = $deDatetime
$ws3.Range("G3").Value = $epochDatetime
$ws3.Range("H3").Value = "IsDependency"
$ws3.Range("I3").Value = "e2e\$mdFile"

$ws3.Range("A4").Value = $mdFile
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = $deXlf3
$ws3.Range("D4").Value = $deDatetime
$ws3.Range("G4").Value = $epochDatetime
$ws3.Range("H4").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "$srcBase/e2e/$pngFile1", "", "", $pngFile1)
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$deHandoffBase/$deXlf1", "", "", $deXlf1)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$srcBase/e2e/$pngFile2", "", "", $pngFile2)
$ws3.Hyperlinks.Add($ws3.Range("C3"), "$deHandoffBase/$deXlf2", "", "", $deXlf2)
$ws3.Hyperlinks.Add($ws3.Range("A4"), "$srcBase/e2e/$mdFile", "", "", $mdFile)
$ws3.Hyperlinks.Add($ws3.Range("C4"), "$deHandoffBase/$deXlf3", "", "", $deXlf3)
$ws3.Hyperlinks.Add($ws3.Range("A5"), "$srcBase/$cfgFile", "", "", $cfgFile)

Write-Host "Localization status report regenerated for handoff."
